# Replace the numeric month values (1-12) in column C (rows 6-85) of the
# "Mes" table column with their Spanish three-letter abbreviation, e.g.
# 8 -> "Ago.", 7 -> "Jul.", ... 1 -> "Ene.", 12 -> "Dic.", 11 -> "Nov.", 9 -> "Sep."
#
# The sequence repeats every 12 rows starting at row 6 (Aug 2024 downward to
# Jan 2024, then Dec 2023 downward, etc.), exactly mirroring the existing
# numeric sequence already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 6; $row -le 85; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C = "Mes"
    $monthNumber = [int]$cell.Value2
    $cell.Value = $monthNames[$monthNumber]
}
